$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("ECs", "Il1b", "Il1r2", "ECs", 3, 1, 1802.141101666667, 5406.423305, 0.8789810559109003, 0.8789810559109003, 2, 0.6666666666666666, 263.9035463333333, 791.710639, 0.9572387917213622, 0.9572387917213622, 475591.4277228936, 4280322.849506042, 0.8413947639061172, 0.8413947639061172)
    ,@("ECs", "Il1b", "Il1r2", "FAPs", 3, 1, 1802.141101666667, 5406.423305, 0.8789810559109003, 0.8789810559109003, 3, 1, 0.9720173333333334, 2.916052, 0.003525730179150291, 0.003525730179150291, 1751.712387932429, 15765.41149139186, 0.003099050035726451, 0.003099050035726451)
    ,@("ECs", "Il1b", "Il1r2", "M2", 3, 1, 1802.141101666667, 5406.423305, 0.8789810559109003, 0.8789810559109003, 3, 1, 9.990833333333333, 29.9725, 0.03623904779290016, 0.03623904779290017, 18004.89138990139, 162044.0225091125, 0.03185343649420896, 0.03185343649420896)
    ,@("ECs", "Il1b", "Il1r2", "sCs", 3, 1, 1802.141101666667, 5406.423305, 0.8789810559109003, 0.8789810559109003, 3, 1, 0.8260933333333332, 2.47828, 0.002996430306587325, 0.002996430306587325, 1488.736749812822, 13398.6307483154, 0.00263380547484755, 0.00263380547484755)
    ,@("M2", "Il1b", "Il1r2", "ECs", 3, 1, 248.1049756666667, 744.314927, 0.121011375461416, 0.121011375461416, 2, 0.6666666666666666, 263.9035463333333, 791.710639, 0.9572387917213622, 0.9572387917213622, 65475.7829413787, 589282.0464724083, 0.1158367828312259, 0.1158367828312259)
    ,@("M2", "Il1b", "Il1r2", "FAPs", 3, 1, 248.1049756666667, 744.314927, 0.121011375461416, 0.121011375461416, 3, 1, 0.9720173333333334, 2.916052, 0.003525730179150291, 0.003525730179150291, 241.1623368342449, 2170.461031508204, 0.0004266534584848014, 0.0004266534584848014)
    ,@("M2", "Il1b", "Il1r2", "M2", 3, 1, 248.1049756666667, 744.314927, 0.121011375461416, 0.121011375461416, 3, 1, 9.990833333333333, 29.9725, 0.03623904779290016, 0.03623904779290017, 2478.775461056389, 22308.9791495075, 0.004385337018830839, 0.004385337018830841)
    ,@("M2", "Il1b", "Il1r2", "sCs", 3, 1, 248.1049756666667, 744.314927, 0.121011375461416, 0.121011375461416, 3, 1, 0.8260933333333332, 2.47828, 0.002996430306587325, 0.002996430306587325, 204.9578663650622, 1844.62079728556, 0.0003626021528744046, 0.0003626021528744047)
    ,@("sCs", "Il1b", "Il1r2", "ECs", 1, 0.3333333333333333, 0.01551766666666667, 0.046553, 0.000007568627683662319, 0.00000756862768366232, 2, 0.6666666666666666, 263.9035463333333, 791.710639, 0.9572387917213622, 0.9572387917213622, 4.095167264151889, 36.856505377367, 0.000007244984018897771, 0.000007244984018897771)
    ,@("sCs", "Il1b", "Il1r2", "FAPs", 1, 0.3333333333333333, 0.01551766666666667, 0.046553, 0.000007568627683662319, 0.00000756862768366232, 3, 1, 0.9720173333333334, 2.916052, 0.003525730179150291, 0.003525730179150291, 0.01508344097288889, 0.135750968756, 0.0000000266849390390406, 0.00000002668493903904061)
    ,@("sCs", "Il1b", "Il1r2", "M2", 1, 0.3333333333333333, 0.01551766666666667, 0.046553, 0.000007568627683662319, 0.00000756862768366232, 3, 1, 9.990833333333333, 29.9725, 0.03623904779290016, 0.03623904779290017, 0.1550344213888889, 1.3953097925, 0.000000274279860354906, 0.0000002742798603549061)
    ,@("sCs", "Il1b", "Il1r2", "sCs", 1, 0.3333333333333333, 0.01551766666666667, 0.046553, 0.000007568627683662319, 0.00000756862768366232, 3, 1, 0.8260933333333332, 2.47828, 0.002996430306587325, 0.002996430306587325, 0.01281904098222222, 0.11537136884, 0.0000000226788653706016, 0.0000000226788653706016)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
